# Updates template Excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "test_semi" to "mal"
$ws.Name = "mal"

# Insert a new (blank) row above the old row 3, pushing the header row and
# all data rows down by one.
$ws.Rows("3").Insert()

# --- Row 1/2: chart axis-title labels -------------------------------------
$ws.Range("A1").Value = "Aksetittel, x-akse"
$ws.Range("A2").Value = "Aksetittel, y-akse"
$ws.Range("B2").Value = "Tonn pr. Måned"

# --- Row 5 (former row 4): table header ------------------------------------
$ws.Range("B5").Value = "Epler"
$ws.Range("C5").Value = "Pærer"
$ws.Range("D5").Value = "Bananer"
$ws.Range("A5").Formula = "=B1"

# --- Row 3 (new): source note -----------------------------------------------
$ws.Range("A3").Value = "Kilde"
$ws.Range("B3").Value = "NDLA Seksjon for eksempelstatistikk"

# --- Updated data values -----------------------------------------------------
# Mars (row 8): D value -1500.5 -> 2000
$ws.Range("D8").Value = 2000
# April (row 9): B value 666 -> 777, C value 2345 -> 1250
$ws.Range("B9").Value = 777
$ws.Range("C9").Value = 1250

# --- View / selection ---------------------------------------------------------
$ws.Range("D28").Select() | Out-Null

# Column A was widened (best-fit) to accommodate the new, longer labels.
$ws.Columns("A").AutoFit()
